# Updated cryptos list on Fri May 17 22:34:26 UTC 2024 with GitHub Actions
#
# Refreshes the "Price" (D) and "Volume(1h)" (E) columns of the crypto
# ticker sheet with newly-scraped values. Both columns hold plain text in
# the source data (e.g. thousand-separator-dotted prices like "66.730.17"
# and padded percentage strings like "  +2.27%  "), so every write below
# is a literal text replacement.
#
# Some new "Price" values happen to be syntactically valid numbers (e.g.
# "580.04"); Excel's COM layer would otherwise auto-coerce those into the
# Number type on assignment. To keep them as text (matching the original
# cell type), we momentarily switch the cell to the Text number format,
# assign the value, then clear the format back so the cell's appearance/
# style is left exactly as it was.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.730.17"
$ws.Range("E2").Value = "  +2.27%  "

$ws.Range("D3").Value = "3.091.17"
$ws.Range("E3").Value = "  +5.51%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.04"
$ws.Range("D5").ClearFormats()

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.03"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +6.48%  "

$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("D8").Value = "3.087.57"
$ws.Range("E8").Value = "  +5.59%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.524"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.90%  "

$ws.Range("E10").Value = "  -1.69%  "

$ws.Range("E11").Value = "  +3.76%  "

$ws.Range("E12").Value = "  +4.83%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000250"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.78%  "

$ws.Range("E14").Value = "  +6.40%  "

$ws.Range("D16").Value = "3.603.39"
$ws.Range("E16").Value = "  +5.43%  "

$ws.Range("D17").Value = "66.718.11"
$ws.Range("E17").Value = "  +2.27%  "

$ws.Range("E18").Value = "  +3.09%  "

$ws.Range("D19").Value = "3.091.39"
$ws.Range("E19").Value = "  +5.57%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.25"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +4.43%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "466.83"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +5.42%  "

$ws.Range("E22").Value = "  +3.56%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.51"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +3.63%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.94"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.27%  "

$ws.Range("E25").Value = "  +6.06%  "

$ws.Range("E26").Value = "  +8.09%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.09"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.75%  "

$ws.Range("E28").Value = "  -0.02%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.01"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.05%  "

$ws.Range("E30").Value = "  +2.45%  "

$ws.Range("E31").Value = "  +4.09%  "

$ws.Range("E32").Value = "  +1.28%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "28.25"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +4.50%  "

$ws.Range("E34").Value = "  +3.81%  "

$ws.Range("E35").Value = "  +0.03%  "

$ws.Range("E36").Value = "  +3.63%  "

$ws.Range("E37").Value = "  +3.08%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "47.37"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +6.00%  "

$ws.Range("E39").Value = "  +6.46%  "

$ws.Range("E40").Value = "  +6.55%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "50.32"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.52%  "

$ws.Range("E42").Value = "  +1.63%  "

$ws.Range("E43").Value = "  +2.49%  "

$ws.Range("E44").Value = "  -0.42%  "

$ws.Range("E45").Value = "  +3.23%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "382.71"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.31%  "

$ws.Range("D47").Value = "2.784.99"
$ws.Range("E47").Value = "  +3.31%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "134.96"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.83%  "

$ws.Range("E49").Value = "  +0.05%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.92"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +7.00%  "

$ws.Range("E51").Value = "  +1.40%  "
